$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed by Excel as a number
# (single-decimal-point numeric strings) are force-formatted as Text first,
# so the stored value exactly matches the source string (e.g. "0.250", "0.500").

$ws.Range('D2').Value = '26.255.01'
$ws.Range('D3').Value = '1.606.37'
$ws.Range('E3').Value = '  -0.17%  '
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.67'
$ws.Range('E5').Value = '  -0.09%  '
$ws.Range('E6').Value = '  -0.15%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.485'
$ws.Range('E7').Value = '  +0.47%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.250'
$ws.Range('E8').Value = '  -0.12%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0615'
$ws.Range('E9').Value = '  -0.62%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '18.23'
$ws.Range('E10').Value = '  +0.47%  '
$ws.Range('D12').Value = '1.826.25'
$ws.Range('E12').Value = '  -0.32%  '
$ws.Range('D13').Value = '1.600.14'
$ws.Range('E13').Value = '  -0.29%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.04'
$ws.Range('E14').Value = '  +0.63%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.514'
$ws.Range('E15').Value = '  +0.93%  '
$ws.Range('D16').Value = '26.264.85'
$ws.Range('E16').Value = '  +0.11%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.53'
$ws.Range('E17').Value = '  +1.41%  '
$ws.Range('D18').Value = '0.0₃0729'
$ws.Range('E18').Value = '  +0.01%  '
$ws.Range('E19').Value = '  -0.15%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '203.19'
$ws.Range('E20').Value = '  +2.07%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.29'
$ws.Range('E21').Value = '  +1.24%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.28'
$ws.Range('E22').Value = '  -1.32%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.99'
$ws.Range('E23').Value = '  -0.54%  '
$ws.Range('E24').Value = '  +9.12%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '144.38'
$ws.Range('E25').Value = '  +1.34%  '
$ws.Range('E26').Value = '  -0.23%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.122'
$ws.Range('E27').Value = '  -6.22%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.21'
$ws.Range('E28').Value = '  +0.09%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.57'
$ws.Range('E29').Value = '  +1.34%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0488'
$ws.Range('E30').Value = '  +3.50%  '
$ws.Range('E31').Value = '  -0.22%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.20'
$ws.Range('E32').Value = '  +2.19%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.94'
$ws.Range('E33').Value = '  -2.49%  '
$ws.Range('E34').Value = '  +2.95%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.49'
$ws.Range('E35').Value = '  -0.04%  '
$ws.Range('D36').Value = '1.151.94'
$ws.Range('E36').Value = '  +3.93%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0165'
$ws.Range('E37').Value = '  +8.44%  '
$ws.Range('E38').Value = '  -0.11%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.796'
$ws.Range('E39').Value = '  +1.18%  '
$ws.Range('E40').Value = '  -0.20%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.500'
$ws.Range('E41').Value = '  -0.18%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.784'
$ws.Range('E42').Value = '  +0.18%  '
$ws.Range('E43').Value = '  +2.76%  '
$ws.Range('D44').Value = '1.739.85'
$ws.Range('E44').Value = '  -0.30%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '91.99'
$ws.Range('E45').Value = '  -0.69%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.53'
$ws.Range('E46').Value = '  -1.68%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '54.20'
$ws.Range('E47').Value = '  +0.84%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0506'
$ws.Range('E48').Value = '  -0.55%  '
$ws.Range('D49').Value = '0.0₇0975'
$ws.Range('E49').Value = '  -9.91%  '
$ws.Range('E50').Value = '  -0.60%  '
$ws.Range('E51').Value = '  -0.27%  '
